$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.834.56"
$ws.Range("E2").Value = "'  -0.08%  "
$ws.Range("D3").Value = "'2.572.43"
$ws.Range("E3").Value = "'  +1.33%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'312.72"
$ws.Range("E5").Value = "'  -0.84%  "
$ws.Range("D6").Value = "'98.62"
$ws.Range("E6").Value = "'  +2.52%  "
$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "'  -0.46%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "'  -0.18%  "
$ws.Range("D10").Value = "'35.75"
$ws.Range("E10").Value = "'  -0.90%  "
$ws.Range("D11").Value = "'0.0811"
$ws.Range("E11").Value = "'  +0.20%  "
$ws.Range("D12").Value = "'7.44"
$ws.Range("E12").Value = "'  -2.02%  "
$ws.Range("D13").Value = "'2.967.46"
$ws.Range("E13").Value = "'  +1.44%  "
$ws.Range("E14").Value = "'  -1.34%  "
$ws.Range("E15").Value = "'  +4.57%  "
$ws.Range("D16").Value = "'2.601.92"
$ws.Range("E16").Value = "'  +2.18%  "
$ws.Range("D17").Value = "'0.845"
$ws.Range("E17").Value = "'  -0.88%  "
$ws.Range("D18").Value = "'42.879.45"
$ws.Range("E18").Value = "'  -0.11%  "
$ws.Range("D19").Value = "'6.74"
$ws.Range("E19").Value = "'  -1.51%  "
$ws.Range("D20").Value = "'12.50"
$ws.Range("E20").Value = "'  -4.28%  "
$ws.Range("D21").Value = "'0.0₃0962"
$ws.Range("E21").Value = "'  -0.48%  "
$ws.Range("D22").Value = "'69.79"
$ws.Range("D23").Value = "'249.43"
$ws.Range("E23").Value = "'  -1.56%  "
$ws.Range("E24").Value = "'  -0.30%  "
$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = "'  -0.56%  "
$ws.Range("D26").Value = "'27.12"
$ws.Range("E26").Value = "'  +1.71%  "
$ws.Range("E27").Value = "'  -0.02%  "
$ws.Range("E28").Value = "'  -1.13%  "
$ws.Range("D29").Value = "'39.69"
$ws.Range("E29").Value = "'  -3.11%  "
$ws.Range("D30").Value = "'10.25"
$ws.Range("E30").Value = "'  -1.87%  "
$ws.Range("D31").Value = "'158.79"
$ws.Range("E31").Value = "'  +0.84%  "
$ws.Range("D32").Value = "'5.79"
$ws.Range("E32").Value = "'  -2.37%  "
$ws.Range("D33").Value = "'3.34"
$ws.Range("E33").Value = "'  +0.06%  "
$ws.Range("D34").Value = "'2.11"
$ws.Range("E34").Value = "'  -2.63%  "
$ws.Range("D35").Value = "'0.0798"
$ws.Range("E35").Value = "'  +1.90%  "
$ws.Range("D36").Value = "'2.68"
$ws.Range("E36").Value = "'  +0.23%  "
$ws.Range("D37").Value = "'18.58"
$ws.Range("E37").Value = "'  -2.31%  "
$ws.Range("D38").Value = "'2.58"
$ws.Range("E38").Value = "'  +10.81%  "
$ws.Range("E39").Value = "'  -0.53%  "
$ws.Range("E40").Value = "'  -0.57%  "
$ws.Range("D41").Value = "'23.02"
$ws.Range("E41").Value = "'  -0.05%  "
$ws.Range("D42").Value = "'4.12"
$ws.Range("E42").Value = "'  +6.85%  "
$ws.Range("B43").Value = "'FirstDigitalUSD"
$ws.Range("C43").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "'  -0.07%  "
$ws.Range("B44").Value = "'VeChain"
$ws.Range("C44").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0302"
$ws.Range("E44").Value = "'  -0.79%  "
$ws.Range("D45").Value = "'3.23"
$ws.Range("E45").Value = "'  -2.43%  "
$ws.Range("D46").Value = "'1.999.11"
$ws.Range("E46").Value = "'  -1.72%  "
$ws.Range("D47").Value = "'9.03"
$ws.Range("E47").Value = "'  -1.56%  "
$ws.Range("D48").Value = "'2.818.64"
$ws.Range("E48").Value = "'  +1.44%  "
$ws.Range("D49").Value = "'0.196"
$ws.Range("E49").Value = "'  +2.02%  "
$ws.Range("D50").Value = "'81.72"
$ws.Range("E50").Value = "'  -4.07%  "
$ws.Range("D51").Value = "'74.40"
$ws.Range("E51").Value = "'  -0.32%  "
